$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
try {
    $cbs = $ws.CheckBoxes()
    Write-Host "CheckBoxes Count: $($cbs.Count)"
} catch {
    Write-Host "Error: $_"
}
